# Auto-generated edit script: update cryptos price/volume table
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values look like plain numbers ("228.89", "1.00", ...).
# Excel auto-converts such text to a numeric value on assignment, but the
# source data must stay a text string (matches the other Price-column cells,
# e.g. "38.797.34", which are never real numbers). Force text format first,
# then restore the default "Normal" style so no stray number format sticks
# to the cell once the literal text value is safely in place.
$textForceCells = @("D5", "D6", "D7", "D10", "D12", "D14", "D15", "D16", "D19", "D20", "D22", "D25", "D26", "D30", "D31", "D33", "D35", "D36", "D37", "D38", "D39", "D40", "D41", "D42", "D44", "D47", "D48")
foreach ($ref in $textForceCells) {
    $ws.Range($ref).NumberFormat = "@"
}

# Apply the updated Price / Volume(1h) values (and the Aave/Maker row swap).
$ws.Range("D2").Value = "38.797.34"
$ws.Range("E2").Value = "  +1.54%  "
$ws.Range("D3").Value = "2.093.27"
$ws.Range("E3").Value = "  -0.07%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").Value = "228.89"
$ws.Range("E5").Value = "  -0.35%  "
$ws.Range("D6").Value = "0.617"
$ws.Range("D7").Value = "61.55"
$ws.Range("E7").Value = "  +1.11%  "
$ws.Range("E8").Value = "  -0.06%  "
$ws.Range("E9").Value = "  +1.57%  "
$ws.Range("D10").Value = "0.0845"
$ws.Range("E10").Value = "  +0.17%  "
$ws.Range("E11").Value = "  -0.15%  "
$ws.Range("D12").Value = "15.31"
$ws.Range("E12").Value = "  +4.37%  "
$ws.Range("D13").Value = "2.403.38"
$ws.Range("E13").Value = "  -0.14%  "
$ws.Range("D14").Value = "22.07"
$ws.Range("D15").Value = "0.806"
$ws.Range("E15").Value = "  +4.06%  "
$ws.Range("D16").Value = "5.48"
$ws.Range("E16").Value = "  -0.78%  "
$ws.Range("D17").Value = "2.090.00"
$ws.Range("E17").Value = "  -0.75%  "
$ws.Range("D18").Value = "38.729.51"
$ws.Range("E18").Value = "  +1.58%  "
$ws.Range("D19").Value = "71.92"
$ws.Range("E19").Value = "  +2.33%  "
$ws.Range("D20").Value = "6.07"
$ws.Range("E20").Value = "  +1.23%  "
$ws.Range("D21").Value = "0.0₃0841"
$ws.Range("E21").Value = "  +0.64%  "
$ws.Range("D22").Value = "228.10"
$ws.Range("E22").Value = "  +1.64%  "
$ws.Range("E23").Value = "  -0.44%  "
$ws.Range("E24").Value = "  -2.57%  "
$ws.Range("D25").Value = "2.34"
$ws.Range("E25").Value = "  +0.46%  "
$ws.Range("D26").Value = "171.47"
$ws.Range("E26").Value = "  +0.69%  "
$ws.Range("E27").Value = "  +0.91%  "
$ws.Range("E28").Value = "  +4.76%  "
$ws.Range("E29").Value = "  +5.78%  "
$ws.Range("D30").Value = "19.31"
$ws.Range("D31").Value = "2.46"
$ws.Range("E31").Value = "  +3.44%  "
$ws.Range("E32").Value = "  +0.87%  "
$ws.Range("D33").Value = "4.51"
$ws.Range("E33").Value = "  +1.80%  "
$ws.Range("E34").Value = "  +0.96%  "
$ws.Range("D35").Value = "0.0619"
$ws.Range("E35").Value = "  +2.25%  "
$ws.Range("D36").Value = "6.50"
$ws.Range("E36").Value = "  -0.28%  "
$ws.Range("D37").Value = "2.38"
$ws.Range("E37").Value = "  -0.88%  "
$ws.Range("D38").Value = "3.60"
$ws.Range("E38").Value = "  +1.58%  "
$ws.Range("D39").Value = "1.00"
$ws.Range("E39").Value = "  +0.12%  "
$ws.Range("D40").Value = "18.26"
$ws.Range("E40").Value = "  +1.07%  "
$ws.Range("D41").Value = "0.0228"
$ws.Range("E41").Value = "  +4.01%  "
$ws.Range("B42").Value = "Aave"
$ws.Range("C42").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D42").Value = "101.15"
$ws.Range("E42").Value = "  +0.86%  "
$ws.Range("B43").Value = "Maker"
$ws.Range("C43").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D43").Value = "1.534.08"
$ws.Range("E43").Value = "  -1.05%  "
$ws.Range("D44").Value = "2.81"
$ws.Range("E44").Value = "  -1.08%  "
$ws.Range("E45").Value = "  +0.45%  "
$ws.Range("E46").Value = "  +1.94%  "
$ws.Range("D47").Value = "7.65"
$ws.Range("E47").Value = "  +5.50%  "
$ws.Range("D48").Value = "4.10"
$ws.Range("E48").Value = "  -0.42%  "
$ws.Range("E49").Value = "  +1.12%  "
$ws.Range("E50").Value = "  -1.12%  "
$ws.Range("D51").Value = "2.290.16"
$ws.Range("E51").Value = "  -0.09%  "

# Drop back to the default style now that the text values are committed,
# so the cells end up with no explicit style index (same as before the edit).
foreach ($ref in $textForceCells) {
    $ws.Range($ref).Style = "Normal"
}
